$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.542630672454834
$ws.Range("B1").Value = 5.143352508544922
$ws.Range("C1").Value = 8.193857192993164
$ws.Range("D1").Value = 7.643161296844482
$ws.Range("E1").Value = 3.78064227104187
